# Update scripts with new TPM (transcripts-per-million) values.
# NATMI LR-pairs output for Efna3-Ephb1: the underlying expression table was
# regenerated, shrinking the cluster-pair list from 6 rows to 4 (the two
# MuSCs -> MuSCs self-pairs at the bottom are gone) and refreshing every
# numeric column for the remaining rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the last two data rows (former rows 6 and 7) ---------------------
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()

# --- Row 2: ECs -> Efna3/Ephb1 -> ECs ---------------------------------------
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 0.3193606666666667
$ws.Range("H2").Value2 = 0.958082
$ws.Range("I2").Value2 = 0.7979421849584948
$ws.Range("J2").Value2 = 0.7979421849584948
$ws.Range("M2").Value2 = 2.718682666666667
$ws.Range("N2").Value2 = 8.156048
$ws.Range("O2").Value2 = 0.5434637507613679
$ws.Range("P2").Value2 = 0.5434637507613679
$ws.Range("Q2").Value2 = 0.8682403088817778
$ws.Range("R2").Value2 = 7.814162779936
$ws.Range("S2").Value2 = 0.4336526527282647
$ws.Range("T2").Value2 = 0.4336526527282647

# --- Row 3: ECs -> Efna3/Ephb1 -> MuSCs -------------------------------------
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 0.3193606666666667
$ws.Range("H3").Value2 = 0.958082
$ws.Range("I3").Value2 = 0.7979421849584948
$ws.Range("J3").Value2 = 0.7979421849584948
$ws.Range("M3").Value2 = 2.283827
$ws.Range("N3").Value2 = 6.851481
$ws.Range("O3").Value2 = 0.4565362492386322
$ws.Range("P3").Value2 = 0.4565362492386321
$ws.Range("Q3").Value2 = 0.7293645132713333
$ws.Range("R3").Value2 = 6.564280619442
$ws.Range("S3").Value2 = 0.3642895322302301
$ws.Range("T3").Value2 = 0.3642895322302301

# --- Row 4: FAPs -> Efna3/Ephb1 -> ECs --------------------------------------
$ws.Range("G4").Value2 = 0.08086966666666666
$ws.Range("I4").Value2 = 0.2020578150415052
$ws.Range("J4").Value2 = 0.2020578150415053
$ws.Range("M4").Value2 = 2.718682666666667
$ws.Range("N4").Value2 = 8.156048
$ws.Range("O4").Value2 = 0.5434637507613679
$ws.Range("P4").Value2 = 0.5434637507613679
$ws.Range("Q4").Value2 = 0.2198589610257778
$ws.Range("R4").Value2 = 1.978730649232
$ws.Range("S4").Value2 = 0.1098110980331032
$ws.Range("T4").Value2 = 0.1098110980331032

# --- Row 5: FAPs -> Efna3/Ephb1 -> MuSCs ------------------------------------
$ws.Range("G5").Value2 = 0.08086966666666666
$ws.Range("I5").Value2 = 0.2020578150415052
$ws.Range("J5").Value2 = 0.2020578150415053
$ws.Range("M5").Value2 = 2.283827
$ws.Range("N5").Value2 = 6.851481
$ws.Range("O5").Value2 = 0.4565362492386322
$ws.Range("P5").Value2 = 0.4565362492386321
$ws.Range("Q5").Value2 = 0.1846923282143333
$ws.Range("R5").Value2 = 1.662230953929
$ws.Range("S5").Value2 = 0.09224671700840208
$ws.Range("T5").Value2 = 0.09224671700840208
